$d = $word.ActiveDocument

# Fix 1: "ont particuierem" -> "ont particulierem" (typo fix, insert missing 'l')
$d.Content.Find.Execute("particuierem", $false, $false, $false, $false, $false,
                         $true, 1, $false, "particulierem", 2)

# Fix 2: " pour la" -> " pour le"
$d.Content.Find.Execute(" pour la", $false, $false, $false, $false, $false,
                         $true, 1, $false, " pour le", 2)
